$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting from the existing "Unnamed: 27" header (AB1) onto
# the three new header cells so they match the bold/centered/bordered style.
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)  # xlPasteFormats

# New header labels
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Fill in the team record (same values for every player row, 2 through 40)
for ($r = 2; $r -le 40; $r++) {
    $ws.Cells.Item($r, 29).Value = 88   # AC - Wins
    $ws.Cells.Item($r, 30).Value = 74   # AD - Losses
    $ws.Cells.Item($r, 31).Value = 0    # AE - Ties
}
